$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark. It gets re-created further down
#    (right after the newly inserted comma, before "<ten_dv>").
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Locate the 4th occurrence of the "<ngay_thang>" placeholder - the one
#    inside "Điều 4" just before "<ten_dv> chịu trách nhiệm thi hành Quyết
#    định này./." - and insert a "," immediately after it.
# ---------------------------------------------------------------------------
$rng = $d.Content.Duplicate
$rng.Start = 0
$matchIndex = 0
$targetRange = $null
while ($rng.Find.Execute("<ngay_thang>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $matchIndex = $matchIndex + 1
    if ($matchIndex -eq 4) {
        $targetRange = $d.Range($rng.Start, $rng.End)
        break
    }
    $rng.Start = $rng.End
    $rng.End = $d.Content.End
}

if ($null -eq $targetRange) {
    throw "Could not find the 4th <ngay_thang> occurrence"
}

# Sanity check: this occurrence must be the one immediately followed by
# "<ten_dv>" (separated only by a space), confirming we are in "Điều 4".
$checkRange = $d.Range($targetRange.End, [Math]::Min($targetRange.End + 20, $d.Content.End))
if (-not ($checkRange.Text -match "<ten_dv>")) {
    throw "Unexpected context after <ngay_thang>; aborting to avoid editing the wrong spot"
}

$insertPos = $targetRange.End
$commaRange = $d.Range($insertPos, $insertPos)
$commaRange.InsertAfter(",")

# ---------------------------------------------------------------------------
# 3. Re-create the "_GoBack" bookmark right after the comma (and before the
#    following space / "<ten_dv>").
# ---------------------------------------------------------------------------
$bmPos = $insertPos + 1
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))

Write-Output "edit applied"
